$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "(313227928, Aviv  Levi: 1,-10)"
$ws.Range("B1").Value = "(205807308, Sariel  Basis: -2,-8)"
$ws.Range("C1").Value = "(315891549, Raz  Halaby: -10,5)"
$ws.Range("D1").Value = "(326598423, Ron Cohen: -3,8)"
$ws.Range("E1").Value = "(313925141, Elad   Amer: -6,4)"
$ws.Range("F1").Value = "(315060103, Dan  Mshelh: -6,1)"
$ws.Range("G1").Value = "(305487936, Avihai  Kipnis: -4,1)"

$ws.Range("A3").Value = "cost: 344.2600523786747"
$ws.Range("A4").Value = "time: 45.60857891123924"
